$d = $word.ActiveDocument
$d.Content.Find.Execute("Katherina Mohort", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Dr. med. Katherina Mohort", 2)
